$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 17 by copying row 16 (same market/category/etc.), then set its unique values below
$ws.Range("A16:R16").Copy($ws.Range("A17:R17"))

# Update date (D), volume (J), min/max/avg price (K/L/M), price per kg (P) for each data row
$ws.Range("D2").Value2 = 44425
$ws.Range("J2").Value2 = 30
$ws.Range("K2").Value2 = 13000
$ws.Range("L2").Value2 = 13000
$ws.Range("M2").Value2 = 13000
$ws.Range("P2").Value2 = 1300

$ws.Range("D3").Value2 = 44473
$ws.Range("J3").Value2 = 25
$ws.Range("K3").Value2 = 11000
$ws.Range("L3").Value2 = 11000
$ws.Range("M3").Value2 = 11000
$ws.Range("P3").Value2 = 1100

$ws.Range("D4").Value2 = 44663
$ws.Range("J4").Value2 = 30
$ws.Range("K4").Value2 = 12000
$ws.Range("L4").Value2 = 12000
$ws.Range("M4").Value2 = 12000
$ws.Range("P4").Value2 = 1200

$ws.Range("D5").Value2 = 44659
$ws.Range("J5").Value2 = 25
$ws.Range("K5").Value2 = 10000
$ws.Range("L5").Value2 = 10000
$ws.Range("M5").Value2 = 10000
$ws.Range("P5").Value2 = 1000

$ws.Range("D6").Value2 = 44649
$ws.Range("J6").Value2 = 25
$ws.Range("K6").Value2 = 10000
$ws.Range("L6").Value2 = 10000
$ws.Range("M6").Value2 = 10000
$ws.Range("P6").Value2 = 1000

$ws.Range("D7").Value2 = 44525
$ws.Range("J7").Value2 = 20
$ws.Range("K7").Value2 = 9000
$ws.Range("L7").Value2 = 9000
$ws.Range("M7").Value2 = 9000
$ws.Range("P7").Value2 = 900

$ws.Range("D8").Value2 = 44645
$ws.Range("J8").Value2 = 25
$ws.Range("K8").Value2 = 10000
$ws.Range("L8").Value2 = 10000
$ws.Range("M8").Value2 = 10000
$ws.Range("P8").Value2 = 1000

$ws.Range("D9").Value2 = 44526
$ws.Range("J9").Value2 = 25
$ws.Range("K9").Value2 = 9000
$ws.Range("L9").Value2 = 9000
$ws.Range("M9").Value2 = 9000
$ws.Range("P9").Value2 = 900

$ws.Range("D10").Value2 = 44348
$ws.Range("J10").Value2 = 20
$ws.Range("K10").Value2 = 10000
$ws.Range("L10").Value2 = 10000
$ws.Range("M10").Value2 = 10000
$ws.Range("P10").Value2 = 1000

$ws.Range("D11").Value2 = 44530
$ws.Range("J11").Value2 = 30
$ws.Range("K11").Value2 = 10000
$ws.Range("L11").Value2 = 10000
$ws.Range("M11").Value2 = 10000
$ws.Range("P11").Value2 = 1000

$ws.Range("D12").Value2 = 44698
$ws.Range("J12").Value2 = 35
$ws.Range("K12").Value2 = 11000
$ws.Range("L12").Value2 = 11000
$ws.Range("M12").Value2 = 11000
$ws.Range("P12").Value2 = 1100

$ws.Range("D13").Value2 = 44656
$ws.Range("J13").Value2 = 25
$ws.Range("K13").Value2 = 10000
$ws.Range("L13").Value2 = 10000
$ws.Range("M13").Value2 = 10000
$ws.Range("P13").Value2 = 1000

$ws.Range("D14").Value2 = 44469
$ws.Range("J14").Value2 = 20
$ws.Range("K14").Value2 = 12000
$ws.Range("L14").Value2 = 12000
$ws.Range("M14").Value2 = 12000
$ws.Range("P14").Value2 = 1200

$ws.Range("D15").Value2 = 44463
$ws.Range("J15").Value2 = 25
$ws.Range("K15").Value2 = 12000
$ws.Range("L15").Value2 = 12000
$ws.Range("M15").Value2 = 12000
$ws.Range("P15").Value2 = 1200

$ws.Range("D16").Value2 = 44369
$ws.Range("J16").Value2 = 25
$ws.Range("K16").Value2 = 8000
$ws.Range("L16").Value2 = 8000
$ws.Range("M16").Value2 = 8000
$ws.Range("P16").Value2 = 800

$ws.Range("D17").Value2 = 44523
$ws.Range("J17").Value2 = 30
$ws.Range("K17").Value2 = 9000
$ws.Range("L17").Value2 = 9000
$ws.Range("M17").Value2 = 9000
$ws.Range("P17").Value2 = 900
